$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.4918566666666667
$ws.Range("H2").Value = 1.47557
$ws.Range("I2").Value = 0.002340719790093636
$ws.Range("J2").Value = 0.002340719790093637
$ws.Range("M2").Value = 9.278280333333333
$ws.Range("N2").Value = 27.834841
$ws.Range("O2").Value = 0.2122966588143784
$ws.Range("P2").Value = 0.2122966588143784
$ws.Range("Q2").Value = 4.563584037152222
$ws.Range("R2").Value = 41.07225633437
$ws.Range("S2").Value = 0.0004969269906575721
$ws.Range("T2").Value = 0.0004969269906575722

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.4918566666666667
$ws.Range("H3").Value = 1.47557
$ws.Range("I3").Value = 0.002340719790093636
$ws.Range("J3").Value = 0.002340719790093637
$ws.Range("O3").Value = 0.2154323368929792
$ws.Range("P3").Value = 0.2154323368929792
$ws.Range("Q3").Value = 4.63098938637
$ws.Range("R3").Value = 41.67890447733
$ws.Range("S3").Value = 0.0005042667343915158
$ws.Range("T3").Value = 0.0005042667343915158

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.4918566666666667
$ws.Range("H4").Value = 1.47557
$ws.Range("I4").Value = 0.002340719790093636
$ws.Range("J4").Value = 0.002340719790093637
$ws.Range("M4").Value = 7.033255
$ws.Range("N4").Value = 21.099765
$ws.Range("O4").Value = 0.1609281551588013
$ws.Range("P4").Value = 0.1609281551588013
$ws.Range("Q4").Value = 3.459353360116667
$ws.Range("R4").Value = 31.13418024105
$ws.Range("S4").Value = 0.0003766877175634655
$ws.Range("T4").Value = 0.0003766877175634655

$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4918566666666667
$ws.Range("H5").Value = 1.47557
$ws.Range("I5").Value = 0.002340719790093636
$ws.Range("J5").Value = 0.002340719790093637
$ws.Range("M5").Value = 17.977458
$ws.Range("N5").Value = 53.932374
$ws.Range("O5").Value = 0.4113428491338411
$ws.Range("P5").Value = 0.411342849133841
$ws.Range("Q5").Value = 8.84233256702
$ws.Range("R5").Value = 79.58099310317999
$ws.Range("S5").Value = 0.0009628383474810827
$ws.Range("T5").Value = 0.0009628383474810827

$ws.Range("I6").Value = 0.03187210184095739
$ws.Range("J6").Value = 0.03187210184095739
$ws.Range("M6").Value = 9.278280333333333
$ws.Range("N6").Value = 27.834841
$ws.Range("O6").Value = 0.2122966588143784
$ws.Range("P6").Value = 0.2122966588143784
$ws.Range("Q6").Value = 62.13943924747377
$ws.Range("R6").Value = 559.254953227264
$ws.Range("S6").Value = 0.006766340730226853
$ws.Range("T6").Value = 0.006766340730226853

$ws.Range("I7").Value = 0.03187210184095739
$ws.Range("J7").Value = 0.03187210184095739
$ws.Range("O7").Value = 0.2154323368929792
$ws.Range("P7").Value = 0.2154323368929792
$ws.Range("S7").Value = 0.006866281381288474
$ws.Range("T7").Value = 0.006866281381288473

$ws.Range("I8").Value = 0.03187210184095739
$ws.Range("J8").Value = 0.03187210184095739
$ws.Range("M8").Value = 7.033255
$ws.Range("N8").Value = 21.099765
$ws.Range("O8").Value = 0.1609281551588013
$ws.Range("P8").Value = 0.1609281551588013
$ws.Range("Q8").Value = 47.10382808917333
$ws.Range("R8").Value = 423.93445280256
$ws.Range("S8").Value = 0.005129118550298707
$ws.Range("T8").Value = 0.005129118550298706

$ws.Range("I9").Value = 0.03187210184095739
$ws.Range("J9").Value = 0.03187210184095739
$ws.Range("M9").Value = 17.977458
$ws.Range("N9").Value = 53.932374
$ws.Range("O9").Value = 0.4113428491338411
$ws.Range("P9").Value = 0.411342849133841
$ws.Range("Q9").Value = 120.400453433344
$ws.Range("R9").Value = 1083.604080900096
$ws.Range("S9").Value = 0.01311036117914335
$ws.Range("T9").Value = 0.01311036117914335

$ws.Range("G10").Value = 6.424796000000001
$ws.Range("H10").Value = 19.274388
$ws.Range("I10").Value = 0.03057526341247335
$ws.Range("J10").Value = 0.03057526341247335
$ws.Range("M10").Value = 9.278280333333333
$ws.Range("N10").Value = 27.834841
$ws.Range("O10").Value = 0.2122966588143784
$ws.Range("P10").Value = 0.2122966588143784
$ws.Range("Q10").Value = 59.61105837247867
$ws.Range("R10").Value = 536.4995253523081
$ws.Range("S10").Value = 0.006491026264837603
$ws.Range("T10").Value = 0.006491026264837603

$ws.Range("G11").Value = 6.424796000000001
$ws.Range("H11").Value = 19.274388
$ws.Range("I11").Value = 0.03057526341247335
$ws.Range("J11").Value = 0.03057526341247335
$ws.Range("O11").Value = 0.2154323368929792
$ws.Range("P11").Value = 0.2154323368929792
$ws.Range("Q11").Value = 60.491529549108
$ws.Range("R11").Value = 544.4237659419721
$ws.Range("S11").Value = 0.00658690044806754
$ws.Range("T11").Value = 0.006586900448067539

$ws.Range("G12").Value = 6.424796000000001
$ws.Range("H12").Value = 19.274388
$ws.Range("I12").Value = 0.03057526341247335
$ws.Range("J12").Value = 0.03057526341247335
$ws.Range("M12").Value = 7.033255
$ws.Range("N12").Value = 21.099765
$ws.Range("O12").Value = 0.1609281551588013
$ws.Range("P12").Value = 0.1609281551588013
$ws.Range("Q12").Value = 45.18722859098001
$ws.Range("R12").Value = 406.68505731882
$ws.Range("S12").Value = 0.004920420734463732
$ws.Range("T12").Value = 0.004920420734463731

$ws.Range("G13").Value = 6.424796000000001
$ws.Range("H13").Value = 19.274388
$ws.Range("I13").Value = 0.03057526341247335
$ws.Range("J13").Value = 0.03057526341247335
$ws.Range("M13").Value = 17.977458
$ws.Range("N13").Value = 53.932374
$ws.Range("O13").Value = 0.4113428491338411
$ws.Range("P13").Value = 0.411342849133841
$ws.Range("Q13").Value = 115.501500248568
$ws.Range("R13").Value = 1039.513502237112
$ws.Range("S13").Value = 0.01257691596510448
$ws.Range("T13").Value = 0.01257691596510447

$ws.Range("G14").Value = 196.516566
$ws.Range("H14").Value = 589.549698
$ws.Range("I14").Value = 0.9352119149564756
$ws.Range("J14").Value = 0.9352119149564756
$ws.Range("M14").Value = 9.278280333333333
$ws.Range("N14").Value = 27.834841
$ws.Range("O14").Value = 0.2122966588143784
$ws.Range("P14").Value = 0.2122966588143784
$ws.Range("Q14").Value = 1823.335789492002
$ws.Range("R14").Value = 16410.02210542802
$ws.Range("S14").Value = 0.1985423648286564
$ws.Range("T14").Value = 0.1985423648286564

$ws.Range("G15").Value = 196.516566
$ws.Range("H15").Value = 589.549698
$ws.Range("I15").Value = 0.9352119149564756
$ws.Range("J15").Value = 0.9352119149564756
$ws.Range("O15").Value = 0.2154323368929792
$ws.Range("P15").Value = 0.2154323368929792
$ws.Range("Q15").Value = 1850.266943740818
$ws.Range("R15").Value = 16652.40249366736
$ws.Range("S15").Value = 0.2014748883292317
$ws.Range("T15").Value = 0.2014748883292316

$ws.Range("G16").Value = 196.516566
$ws.Range("H16").Value = 589.549698
$ws.Range("I16").Value = 0.9352119149564756
$ws.Range("J16").Value = 0.9352119149564756
$ws.Range("M16").Value = 7.033255
$ws.Range("N16").Value = 21.099765
$ws.Range("O16").Value = 0.1609281551588013
$ws.Range("P16").Value = 0.1609281551588013
$ws.Range("Q16").Value = 1382.15112040233
$ws.Range("R16").Value = 12439.36008362097
$ws.Range("S16").Value = 0.1505019281564754
$ws.Range("T16").Value = 0.1505019281564754

$ws.Range("G17").Value = 196.516566
$ws.Range("H17").Value = 589.549698
$ws.Range("I17").Value = 0.9352119149564756
$ws.Range("J17").Value = 0.9352119149564756
$ws.Range("M17").Value = 17.977458
$ws.Range("N17").Value = 53.932374
$ws.Range("O17").Value = 0.4113428491338411
$ws.Range("P17").Value = 0.411342849133841
$ws.Range("Q17").Value = 3532.868311569228
$ws.Range("R17").Value = 31795.81480412305
$ws.Range("S17").Value = 0.3846927336421122
$ws.Range("T17").Value = 0.3846927336421121
